# Updated capital structure database
# Refreshes the computed capital-structure / profitability metrics for the
# three Mauritius "Financial Svcs. (Non-bank & Insurance)" rows, and drops
# a handful of cells (E3, T2:T4) that are no longer populated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (industry aggregate) ---------------------------------------
$ws.Range("D2").Value  = -0.173
$ws.Range("E2").Value  = -0.198
$ws.Range("K2").Value  = 2.7
$ws.Range("L2").Value  = 0.07433920704845816
$ws.Range("M2").Value  = 0
$ws.Range("N2").Value  = 0
$ws.Range("O2").Value  = 0
$ws.Range("P2").Value  = 0
$ws.Range("Q2").Value  = 0
$ws.Range("R2").Value  = 0
$ws.Range("U2").Value  = 11.6
$ws.Range("V2").Value  = 0.07332490518331226
$ws.Range("W2").Value  = 0.01438061866345588
$ws.Range("X2").Value  = 0.03044506089345466
$ws.Range("Y2").Value  = -0.01606444222999878
$ws.Range("Z2").Value  = 0.07019985310603426
$ws.Range("AB2").Value = 0.02706831419273371
$ws.Range("AC2").Value = -0.02706831419273371
$ws.Range("AD2").Value = 261.27
$ws.Range("AF2").Value = 261.27
$ws.Range("AG2").Value = 249.67
$ws.Range("AH2").Value = 0.6228574153097957
$ws.Range("AI2").Value = 0.5152543041394679
$ws.Range("AJ2").Value = 0.612131316350798
$ws.Range("AK2").Value = 0.5039053827678769
$ws.Range("AM2").Value = -6.56

# --- Row 3 (Fincorp Investment Ltd.) -----------------------------------
$ws.Range("K3").Value  = -2.09
$ws.Range("L3").Value  = -4.976190476190476
$ws.Range("M3").Value  = 0
$ws.Range("N3").Value  = 0
$ws.Range("O3").Value  = 0
$ws.Range("P3").Value  = 0
$ws.Range("Q3").Value  = 0
$ws.Range("R3").Value  = 0
$ws.Range("W3").Value  = -0.01321948134092347
$ws.Range("X3").Value  = 0.02125421650162471
$ws.Range("Y3").Value  = -0.03447369784254818
$ws.Range("Z3").Value  = 0.002516778523489933
$ws.Range("AB3").Value = 0.02238275631447389
$ws.Range("AC3").Value = -0.02238275631447389
$ws.Range("AD3").Value = 9.17
$ws.Range("AF3").Value = 9.17
$ws.Range("AG3").Value = 9.17
$ws.Range("AH3").Value = 0.1598396374411713
$ws.Range("AI3").Value = 0.06130908604666711
$ws.Range("AJ3").Value = 0.1598396374411713
$ws.Range("AK3").Value = 0.06130908604666711
$ws.Range("AM3").Value = -6.56

# --- Row 4 (CIM Financial Services Ltd) --------------------------------
$ws.Range("D4").Value  = -0.173
$ws.Range("E4").Value  = -0.198
$ws.Range("K4").Value  = 4.79
$ws.Range("L4").Value  = 0.1334261838440111
$ws.Range("M4").Value  = 0
$ws.Range("N4").Value  = 0
$ws.Range("O4").Value  = 0
$ws.Range("P4").Value  = 0
$ws.Range("Q4").Value  = 0
$ws.Range("R4").Value  = 0
$ws.Range("U4").Value  = 11.6
$ws.Range("V4").Value  = 0.1054545454545455
$ws.Range("W4").Value  = 0.04198071866783523
$ws.Range("X4").Value  = 0.03963590528528461
$ws.Range("Y4").Value  = 0.002344813382550617
$ws.Range("Z4").Value  = 0.1024251069900143
$ws.Range("AB4").Value = 0.03175387207099353
$ws.Range("AC4").Value = -0.03175387207099353
$ws.Range("AD4").Value = 252.1
$ws.Range("AF4").Value = 252.1
$ws.Range("AG4").Value = 240.5
$ws.Range("AH4").Value = 0.696216514774924
$ws.Range("AI4").Value = 0.7051748251748251
$ws.Range("AJ4").Value = 0.6861626248216833
$ws.Range("AK4").Value = 0.6952876553917318

# --- Cells removed entirely in the refreshed export --------------------
$ws.Range("T2").ClearContents()
$ws.Range("E3").ClearContents()
$ws.Range("T3").ClearContents()
$ws.Range("T4").ClearContents()
